$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.892.29'

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.38%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.633.94'

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -1.07%  '

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.10%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '211.83'

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.96%  '

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.81%  '

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.06%  '

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -1.97%  '

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -3.57%  '

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -0.45%  '

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.79%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.865.19'

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -1.01%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.637.27'

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -0.80%  '

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -0.86%  '

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.28%  '

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.96%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '27.900.48'

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.29%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '229.73'

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -1.51%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0₃0720'

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.55%  '

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -2.68%  '

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.03%  '

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.97%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.32'

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -3.64%  '

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -3.63%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '152.59'

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.15%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.93'

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.24%  '

$ws.Range("B27").Value = 'Stellar'

$ws.Range("C27").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.111'

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.71%  '

$ws.Range("B28").Value = 'EthereumClassic'

$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.60'

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -1.06%  '

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +0.11%  '

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -1.21%  '

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.77%  '

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +0.40%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.397.14'

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -4.04%  '

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -1.65%  '

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -0.52%  '

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +9.02%  '

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +1.60%  '

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +0.10%  '

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -0.62%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.869'

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -2.66%  '

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +0.07%  '

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -0.04%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '66.76'

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -3.92%  '

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +2.61%  '

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.04%  '

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -1.71%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.775.18'

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -0.92%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '87.65'

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -1.41%  '

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.75%  '

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -0.23%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.53'

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -2.61%  '

